# Update periods in 24Tto25TMap from 2022Oct to 2023Oct
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

for ($r = 2; $r -le 68; $r++) {
    $cell = $ws.Cells.Item($r, 9)  # Column I
    if ($cell.Value2 -eq "2022Oct") {
        $cell.Value2 = "2023Oct"
    }
}

# Update the active selection to G8 to match the recorded view state
$ws.Range("G8").Select()
